# DEAN import format update:
#  - Enrollment sheet gains 10 new trailing columns (C:M) with new headers
#    (person_id re-used + 9 brand new field names added to sharedStrings).
#  - Column widths set for the new Enrollment columns.
#  - Selection/active-cell state updated on Person, Course_Section and
#    Enrollment sheets; Course_Section becomes the active tab.

$wb = $excel.ActiveWorkbook

# --- Person sheet: selection moves from B8 to D9 -----------------------
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Range("D9").Select()

# --- Enrollment sheet: new header columns + widths + selection ---------
$wsEnrollment = $wb.Worksheets.Item("Enrollment")

$wsEnrollment.Range("C1").Value = "person_id"
$wsEnrollment.Range("D1").Value = "enrollment_date"
$wsEnrollment.Range("E1").Value = "completion_flag"
$wsEnrollment.Range("F1").Value = "completion_success_flag"
$wsEnrollment.Range("G1").Value = "withdrawal_flag"
$wsEnrollment.Range("H1").Value = "drop_flag"
$wsEnrollment.Range("I1").Value = "enrollment_status_change_date"
$wsEnrollment.Range("J1").Value = "course_grade_final_number"
$wsEnrollment.Range("K1").Value = "course_grade_final_letter"
$wsEnrollment.Range("L1").Value = "course_grade_to_date_number"
$wsEnrollment.Range("M1").Value = "course_grade_to_date_letter"

$wsEnrollment.Columns.Item(4).ColumnWidth = 14.2481481481481
$wsEnrollment.Columns.Item(5).ColumnWidth = 13.9481481481481
$wsEnrollment.Columns.Item(6).ColumnWidth = 21.0111111111111
$wsEnrollment.Columns.Item(7).ColumnWidth = 13.7444444444444
$wsEnrollment.Columns.Item(8).ColumnWidth = 9.21111111111111
$wsEnrollment.Columns.Item(9).ColumnWidth = 26.6555555555556
$wsEnrollment.Columns.Item(10).ColumnWidth = 23.5333333333333
$wsEnrollment.Columns.Item(11).ColumnWidth = 21.5111111111111
$wsEnrollment.Columns.Item(12).ColumnWidth = 26.1518518518519
$wsEnrollment.Columns.Item(13).ColumnWidth = 24.1333333333333

$wsEnrollment.Range("K24").Select()

# --- Course_Section sheet: becomes the active tab, selection A11 -------
$wsCourseSection = $wb.Worksheets.Item("Course_Section")
$wsCourseSection.Activate()
$wsCourseSection.Range("A11").Select()
